# Update pins such that left/right were swapped to make layout easy
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F15").Value = "P15"
$ws.Range("G15").Value = 6

$ws.Range("F16").Value = "P55"
$ws.Range("G16").Value = 3

$ws.Range("F17").Value = "P59"
$ws.Range("G17").Value = 3

$ws.Range("F21").Value = "P65"

$ws.Range("F22").Value = "P14"
$ws.Range("G22").Value = 6

$ws.Range("F25").Value = "P61"
$ws.Range("G25").Value = 3

$ws.Range("F26").Value = "P62"
$ws.Range("G26").Value = 3

$ws.Range("F27").Value = "P67"
$ws.Range("G27").Value = 2

$ws.Range("F28").Value = "P68"
$ws.Range("G28").Value = 2

$ws.Range("F29").Value = "P71"
$ws.Range("G29").Value = 2

$ws.Range("F30").Value = "P72"
$ws.Range("G30").Value = 2

$ws.Range("F32").Value = "P9"
$ws.Range("G32").Value = 7

$ws.Range("F33").Value = "P8"
$ws.Range("G33").Value = 7

$ws.Range("F34").Value = "P5"
$ws.Range("G34").Value = 7

$ws.Range("F35").Value = "P4"
$ws.Range("G35").Value = 7

$ws.Range("F39").Value = "P21"
$ws.Range("G39").Value = 6

# Update the view: select F39 (this also resets the scroll position/topLeftCell)
$ws.Activate()
$ws.Range("F39").Select()
